# Update the ultrasound equipment brand/model referenced in the report
# template: "MEDISONIC MODELO H60" -> "MINDRAY MODELO DC – N3"
$d = $word.ActiveDocument
$d.Content.Find.Execute("MEDISONIC MODELO H60 ", $false, $false, $false, $false, $false, $true, 1, $false, "MINDRAY MODELO DC – N3 ", 2)
